$d = $word.ActiveDocument

# --- 1. Title paragraph: "Reducing Gun Violence" -> insert a new empty
#        centered paragraph before it, then retitle the original text to
#        "How to reduce gun violence:" ---------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`n") -eq "Reducing Gun Violence") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Insert a bare paragraph mark right at the start of the target
    # paragraph; this splits the paragraph in two without leaving a stray
    # empty run behind in either half.
    $insertPoint = $d.Range($target.Range.Start, $target.Range.Start)
    $insertPoint.Text = [char]13

    # The original paragraph (now pushed one slot later) still carries the
    # "Reducing Gun Violence" text; the freshly created one just before it
    # is blank and should pick up the same centering as the title block
    # above it.
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.TrimEnd("`r", "`n") -eq "Reducing Gun Violence") {
            $target = $p
            break
        }
    }
    $newPara = $target.Previous()
    $newPara.Alignment = 1

    $target.Range.Text = "How to reduce gun violence:"
}

# --- 2. "Police Response Times" -> "The effect of higher police response
#        times:" ------------------------------------------------------
$d.Content.Find.Execute("Police Response Times", $true, $false, $false, $false, $false,
                         $true, 1, $false, "The effect of higher police response times:", 2) | Out-Null

# --- 3. Final empty list-item paragraph loses its list/style formatting,
#        becoming a plain empty paragraph --------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.ListFormat.RemoveNumbers() | Out-Null
$lastPara.Style = "Normal"
$lastPara.Range.ParagraphFormat.Reset() | Out-Null
